$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.958.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("E9").Value = "  -2.78%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.645.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.952.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("E26").Value = "  +0.91%  "

$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.408.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.51%  "

$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  +9.77%  "

$ws.Range("E37").Value = "  +1.48%  "

$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("E44").Value = "  +2.83%  "

$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.775.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("E51").Value = "  -0.29%  "
